# chore: update Sheets via scheduled runner
#
# Refreshes cached Market Board pricing figures (currentAveragePrice /
# currentAveragePriceNQ / currentAveragePriceHQ / LevePriceNQ / LevePriceHQ /
# LeveProfitNQ / LeveProfitHQ) for a batch of leves across the ALC, ARM, BSM,
# CRP, CUL, GSM, LTW and WVR sheets, as produced by the scheduled price-sync
# job. Profit cells that are undefined when their corresponding price is 0
# are cleared (ClearContents) rather than left at 0, matching upstream's
# convention of omitting N/A profit figures.

$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 2412.125
$ws.Range("J32").Value = 2582.8333
$ws.Range("L32").Value = 2582.8333
$ws.Range("N32").Value = -3234.8333
$ws.Range("H70").Value = 2777.182
$ws.Range("I70").Value = 3399.8
$ws.Range("J70").Value = 2258.3333
$ws.Range("K70").Value = 10199.4
$ws.Range("L70").Value = 6774.999899999999
$ws.Range("M70").Value = -9929.400000000001
$ws.Range("N70").Value = -7314.999899999999
$ws.Range("H73").Value = 2777.182
$ws.Range("I73").Value = 3399.8
$ws.Range("J73").Value = 2258.3333
$ws.Range("K73").Value = 10199.4
$ws.Range("L73").Value = 6774.999899999999
$ws.Range("M73").Value = -9263.400000000001
$ws.Range("N73").Value = -8646.999899999999
$ws.Range("H118").Value = 4036.4285
$ws.Range("I118").Value = 4036.4285
$ws.Range("J118").Value = 0
$ws.Range("K118").Value = 12109.2855
$ws.Range("L118").Value = 0
$ws.Range("M118").Value = -10452.2855
$ws.Range("N118").ClearContents()
$ws.Range("H132").Value = 2524.3103
$ws.Range("I132").Value = 2611.7307
$ws.Range("K132").Value = 7835.1921
$ws.Range("M132").Value = -5305.1921
$ws.Range("H135").Value = 1249.2354
$ws.Range("I135").Value = 518.2308
$ws.Range("K135").Value = 4664.077200000001
$ws.Range("M135").Value = -2129.077200000001
$ws.Range("H141").Value = 2683
$ws.Range("I141").Value = 1983.75
$ws.Range("K141").Value = 5951.25
$ws.Range("M141").Value = -771.25

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 5211.9697
$ws.Range("I61").Value = 5158.625
$ws.Range("K61").Value = 5158.625
$ws.Range("M61").Value = -4946.625
$ws.Range("H107").Value = 199500
$ws.Range("J107").Value = 199500
$ws.Range("L107").Value = 199500
$ws.Range("N107").Value = -207180
$ws.Range("H109").Value = 199500
$ws.Range("J109").Value = 199500
$ws.Range("L109").Value = 199500
$ws.Range("N109").Value = -202274
$ws.Range("H122").Value = 2985.7666
$ws.Range("I122").Value = 1890.2174
$ws.Range("K122").Value = 5670.6522
$ws.Range("M122").Value = -3220.6522
$ws.Range("H132").Value = 4029.6553
$ws.Range("I132").Value = 2457.4092
$ws.Range("J132").Value = 8971
$ws.Range("K132").Value = 7372.2276
$ws.Range("L132").Value = 26913
$ws.Range("M132").Value = -4842.2276
$ws.Range("N132").Value = -31973
$ws.Range("H136").Value = 5211.9697
$ws.Range("I136").Value = 5158.625
$ws.Range("K136").Value = 15475.875
$ws.Range("M136").Value = -12925.875

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H29").Value = 3255
$ws.Range("I29").Value = 3255
$ws.Range("K29").Value = 3255
$ws.Range("M29").Value = -2966
$ws.Range("H59").Value = 48663.332
$ws.Range("J59").Value = 62995
$ws.Range("L59").Value = 62995
$ws.Range("N59").Value = -64689
$ws.Range("H108").Value = 182916.67
$ws.Range("J108").Value = 182916.67
$ws.Range("L108").Value = 182916.67
$ws.Range("N108").Value = -190596.67

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1490.7693
$ws.Range("I16").Value = 902.6
$ws.Range("K16").Value = 902.6
$ws.Range("M16").Value = -615.6
$ws.Range("H31").Value = 6498.0938
$ws.Range("I31").Value = 2725.8572
$ws.Range("K31").Value = 2725.8572
$ws.Range("M31").Value = -2430.8572
$ws.Range("H34").Value = 6498.0938
$ws.Range("I34").Value = 2725.8572
$ws.Range("K34").Value = 2725.8572
$ws.Range("M34").Value = -2523.8572
$ws.Range("H58").Value = 3468.611
$ws.Range("I58").Value = 1953.7273
$ws.Range("J58").Value = 5849.143
$ws.Range("K58").Value = 1953.7273
$ws.Range("L58").Value = 5849.143
$ws.Range("M58").Value = -1750.7273
$ws.Range("N58").Value = -6255.143
$ws.Range("H62").Value = 4980
$ws.Range("I62").Value = 3966.6667
$ws.Range("J62").Value = 6500
$ws.Range("K62").Value = 3966.6667
$ws.Range("L62").Value = 6500
$ws.Range("M62").Value = -3342.6667
$ws.Range("N62").Value = -7748
$ws.Range("H65").Value = 4980
$ws.Range("I65").Value = 3966.6667
$ws.Range("J65").Value = 6500
$ws.Range("K65").Value = 19833.3335
$ws.Range("L65").Value = 32500
$ws.Range("M65").Value = -16713.3335
$ws.Range("N65").Value = -38740
$ws.Range("H99").Value = 3579.5
$ws.Range("I99").Value = 3536.875
$ws.Range("J99").Value = 3750
$ws.Range("K99").Value = 3536.875
$ws.Range("L99").Value = 3750
$ws.Range("M99").Value = -2038.875
$ws.Range("N99").Value = -6746
$ws.Range("H113").Value = 1490.7693
$ws.Range("I113").Value = 902.6
$ws.Range("K113").Value = 902.6
$ws.Range("M113").Value = 1267.4
$ws.Range("H122").Value = 3441.8635
$ws.Range("I122").Value = 3063.8
$ws.Range("J122").Value = 4252
$ws.Range("K122").Value = 9191.400000000001
$ws.Range("L122").Value = 12756
$ws.Range("M122").Value = -6741.400000000001
$ws.Range("N122").Value = -17656
$ws.Range("H126").Value = 3579.5
$ws.Range("I126").Value = 3536.875
$ws.Range("J126").Value = 3750
$ws.Range("K126").Value = 10610.625
$ws.Range("L126").Value = 11250
$ws.Range("M126").Value = -8140.625
$ws.Range("N126").Value = -16190
$ws.Range("H135").Value = 94974.14
$ws.Range("J135").Value = 94974.14
$ws.Range("L135").Value = 94974.14
$ws.Range("N135").Value = -105114.14
$ws.Range("H136").Value = 3468.611
$ws.Range("I136").Value = 1953.7273
$ws.Range("J136").Value = 5849.143
$ws.Range("K136").Value = 5861.1819
$ws.Range("L136").Value = 17547.429
$ws.Range("M136").Value = -3311.1819
$ws.Range("N136").Value = -22647.429
$ws.Range("H140").Value = 93535.38
$ws.Range("J140").Value = 93535.38
$ws.Range("L140").Value = 93535.38
$ws.Range("N140").Value = -103895.38

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H50").Value = 74
$ws.Range("I50").Value = 36.5
$ws.Range("J50").Value = 149
$ws.Range("K50").Value = 109.5
$ws.Range("L50").Value = 447
$ws.Range("M50").Value = 371.5
$ws.Range("N50").Value = -1409
$ws.Range("H53").Value = 74
$ws.Range("I53").Value = 36.5
$ws.Range("J53").Value = 149
$ws.Range("K53").Value = 109.5
$ws.Range("L53").Value = 447
$ws.Range("M53").Value = 371.5
$ws.Range("N53").Value = -1409
$ws.Range("H61").Value = 106.416664
$ws.Range("I61").Value = 106.416664
$ws.Range("K61").Value = 319.249992
$ws.Range("M61").Value = -104.249992
$ws.Range("H107").Value = 386.6
$ws.Range("J107").Value = 386.6
$ws.Range("L107").Value = 1159.8
$ws.Range("N107").Value = -4999.8

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 5231.0835
$ws.Range("J122").Value = 11649.333
$ws.Range("L122").Value = 34947.999
$ws.Range("N122").Value = -39847.999
$ws.Range("H132").Value = 4098.2964
$ws.Range("I132").Value = 2565.2778
$ws.Range("J132").Value = 7164.3335
$ws.Range("K132").Value = 7695.8334
$ws.Range("L132").Value = 21493.0005
$ws.Range("M132").Value = -5165.8334
$ws.Range("N132").Value = -26553.0005
$ws.Range("H134").Value = 51557.715
$ws.Range("J134").Value = 51557.715
$ws.Range("L134").Value = 154673.145
$ws.Range("N134").Value = -159743.145

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3328.5715
$ws.Range("I46").Value = 909.1
$ws.Range("K46").Value = 909.1
$ws.Range("M46").Value = -721.1
$ws.Range("H69").Value = 155000
$ws.Range("J69").Value = 155000
$ws.Range("L69").Value = 155000
$ws.Range("N69").Value = -156622
$ws.Range("H72").Value = 155000
$ws.Range("J72").Value = 155000
$ws.Range("L72").Value = 465000
$ws.Range("N72").Value = -473112
$ws.Range("H122").Value = 5450.4194
$ws.Range("I122").Value = 4501.76
$ws.Range("K122").Value = 13505.28
$ws.Range("M122").Value = -11055.28
$ws.Range("H132").Value = 5370.273
$ws.Range("I132").Value = 4157.3
$ws.Range("J132").Value = 17500
$ws.Range("K132").Value = 12471.9
$ws.Range("L132").Value = 52500
$ws.Range("M132").Value = -9941.900000000001
$ws.Range("N132").Value = -57560
$ws.Range("H140").Value = 57085.25
$ws.Range("J140").Value = 57085.25
$ws.Range("L140").Value = 57085.25
$ws.Range("N140").Value = -67445.25

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H119").Value = 130408.5
$ws.Range("J119").Value = 161333.33
$ws.Range("L119").Value = 161333.33
$ws.Range("N119").Value = -171009.33
$ws.Range("H122").Value = 2804.5833
$ws.Range("I122").Value = 2875.8667
$ws.Range("J122").Value = 2685.7778
$ws.Range("K122").Value = 8627.6001
$ws.Range("L122").Value = 8057.3334
$ws.Range("M122").Value = -6177.6001
$ws.Range("N122").Value = -12957.3334
$ws.Range("H132").Value = 3300.85
$ws.Range("I132").Value = 2459.8235
$ws.Range("J132").Value = 8066.6665
$ws.Range("K132").Value = 7379.470499999999
$ws.Range("L132").Value = 24199.9995
$ws.Range("M132").Value = -4849.470499999999
$ws.Range("N132").Value = -29259.9995
$ws.Range("H138").Value = 80424
$ws.Range("J138").Value = 80424
$ws.Range("L138").Value = 80424
$ws.Range("N138").Value = -90704
